$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1. Fix up three already-logged time entries (rows 20, 30, 31, 32)
# ------------------------------------------------------------------
$ws1.Range("C20").Value = 0.60416666666666663

$ws1.Range("C30").Value = 0.5625
$ws1.Range("D30").Value = 0.57638888888888895

$ws1.Range("C31").Value = 0.57638888888888895
$ws1.Range("D31").Value = 0.59722222222222221

$ws1.Range("C32").Value = 0.59722222222222221

# ------------------------------------------------------------------
# 2. Grow the journal table (Tableau1) from B3:I35 to B3:I51 and
#    fill in the newly logged "diagramme de sequence" work entries.
# ------------------------------------------------------------------
$lo = $ws1.ListObjects.Item("Tableau1")
$lo.Resize($ws1.Range("B3:I51"))

# Row 35 - close out the previous entry (Analyse)
$ws1.Range("B35").Value = 45415
$ws1.Range("C35").Value = 0.6875
$ws1.Range("D35").Value = 0.70486111111111116
$ws1.Range("F35").Value = "Analyse"

# Row 36 - start of the sequence-diagram work (taller row, wraps text)
$ws1.Rows("36").RowHeight = 30
$ws1.Range("B36").Value = 45418
$ws1.Range("C36").Value = 0.33333333333333331
$ws1.Range("D36").Value = 0.34722222222222227
$ws1.Range("F36").Value = "Analyse"

$ws1.Range("I36").Value = "https://astah.net/support/astah-pro/user-guide/sequence-diagram/"
$ws1.Hyperlinks.Add($ws1.Range("I36"), "https://astah.net/support/astah-pro/user-guide/sequence-diagram/") | Out-Null
$ws1.Range("I17").Copy() | Out-Null
$ws1.Range("I36").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws1.Range("G37").Value = "Création du diagramme de séquence du déplacement d'une pièce"
$ws1.Range("G36").Value = "Création du diagramme de séquence du début de la partie"
$ws1.Range("G38").Value = "Création du diagramme de séquence du fin de jeu"
$ws1.Range("G39").Value = "Réalisation du diagramme de séquence du fin de jeu"
$ws1.Range("G40").Value = "Écriture des diagrammes de séquences dans le dossier de projet"

# Row 37
$ws1.Range("B37").Value = 45418
$ws1.Range("C37").Value = 0.34722222222222227
$ws1.Range("D37").Value = 0.38194444444444442
$ws1.Range("F37").Value = "Analyse"

# Row 38
$ws1.Range("B38").Value = 45418
$ws1.Range("C38").Value = 0.38194444444444442
$ws1.Range("D38").Value = 0.39930555555555558
$ws1.Range("F38").Value = "Analyse"

# Row 39
$ws1.Range("B39").Value = 45418
$ws1.Range("C39").Value = 0.40972222222222227
$ws1.Range("D39").Value = 0.4236111111111111
$ws1.Range("F39").Value = "Analyse"

# Row 40 - documentation step, still in progress (no end time)
$ws1.Range("B40").Value = 45418
$ws1.Range("C40").Value = 0.4236111111111111
$ws1.Range("F40").Value = "Documentation"

# Rows 35-51: make sure the "Duree" calculated column formula covers
# every row now inside the resized table.
for ($r = 35; $r -le 51; $r++) {
    $ws1.Range("E$r").Formula = "=Tableau1[[#This Row],[Fin]]-Tableau1[[#This Row],[Début]]"
}

# ------------------------------------------------------------------
# 3. View-state: the author had switched back to Feuil1 (and
#    scrolled/selected differently on both tabs) before saving.
# ------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("G29").Select()

$ws1.Activate()
$excel.ActiveWindow.Zoom = 100
$ws1.Range("G57").Select()

$wb.Save()
